# Update "想去人数" (want-to-go count) figures in column F across all four sheets,
# matching the latest data refresh.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F16").Value = 92
$ws1.Range("F20").Value = 56942
$ws1.Range("F29").Value = 4551
$ws1.Range("F31").Value = 76
$ws1.Range("F35").Value = 1376
$ws1.Range("F43").Value = 209
$ws1.Range("F47").Value = 47

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 56
$ws2.Range("F31").Value = 1
$ws2.Range("F35").Value = 1
$ws2.Range("F48").Value = 132

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F10").Value = 1711
$ws3.Range("F15").Value = 223
$ws3.Range("F16").Value = 2034
$ws3.Range("F17").Value = 421

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 1711
$ws4.Range("F14").Value = 223
$ws4.Range("F17").Value = 92
$ws4.Range("F19").Value = 56942
$ws4.Range("F27").Value = 76
$ws4.Range("F31").Value = 421
$ws4.Range("F42").Value = 209
$ws4.Range("F47").Value = 47

$wb.Save()
